# Generate Report for Handoff
# Adds two new source files (7169c82a..., 8d2e6537...) to the localization
# status report, ahead of the existing c9099fd2... entry, and refreshes the
# handoff timestamps for all rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Common lookups
# ---------------------------------------------------------------------
$mdBase   = "https://github.com/OpenLocalizationTest/oltest/blob/1d06f000eb2c73d7329417efa25d2835779e1327/e2e/"
$xlfZhBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e390b96cc72575d89428cadd4c54f64a1419d10c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/"
$xlfDeBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/027aa8a974ae759cacd940d9eeb6b6d547b0245c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/"

$file1 = "4a40ba11-58f1-477b-974b-4cc31e6347b4"
$file2 = "7169c82a-aa1b-46f7-94d0-2759d90a762a"
$file3 = "8d2e6537-98be-4a83-befd-cbebfd9a48d1"
$file4 = "c9099fd2-3afe-4f73-b0d0-3a9fb2340893"

$hash1 = "00a2ddbcb0af86db3941c4468cf7db52c709a333"
$hash2 = "2c006fdd7382a1e00a35982b5ac53f371fea4a12"
$hash3 = "1baba9d3437cb8cc2dae41cbe4a76a6fe8228aec"
$hash4 = "64a485ca666abd3c57f37aaf057b78cac10c2d71"

$handoffDate = "2016-12-13 08:12:46"
$zhDatetime  = "2016-03-13 08:12:42"
$deDatetime  = "2016-03-13 08:12:46"
$epoch       = "0001-01-01 00:00:00"

# =======================================================================
# Sheet "Overview"
# =======================================================================
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = $file1 + ".md"
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = $handoffDate

$ws.Range("A3").Value = $file2 + ".md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = $handoffDate

$ws.Range("A4").Value = $file3 + ".md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = $handoffDate

$ws.Range("A5").Value = $file4 + ".md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = $handoffDate

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdBase + $file1 + ".md", "", "", $file1 + ".md")
$ws.Hyperlinks.Add($ws.Range("A3"), $mdBase + $file2 + ".md", "", "", $file2 + ".md")
$ws.Hyperlinks.Add($ws.Range("A4"), $mdBase + $file3 + ".md", "", "", $file3 + ".md")
$ws.Hyperlinks.Add($ws.Range("A5"), $mdBase + $file4 + ".md", "", "", $file4 + ".md")

# =======================================================================
# Sheet "zh-cn"
# =======================================================================
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = $file1 + ".md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = $file1 + "." + $hash1 + ".zh-cn.xlf"
$ws.Range("E2").Value = $zhDatetime
$ws.Range("H2").Value = $epoch
$ws.Range("I2").Value = "Include"

$ws.Range("A3").Value = $file2 + ".md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = $file2 + "." + $hash2 + ".zh-cn.xlf"
$ws.Range("E3").Value = $zhDatetime
$ws.Range("H3").Value = $epoch
$ws.Range("I3").Value = "Include"

$ws.Range("A4").Value = $file3 + ".md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = $file3 + "." + $hash3 + ".zh-cn.xlf"
$ws.Range("E4").Value = $zhDatetime
$ws.Range("H4").Value = $epoch
$ws.Range("I4").Value = "Include"

$ws.Range("A5").Value = $file4 + ".md"
$ws.Range("B5").Value = ".md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = $file4 + "." + $hash4 + ".zh-cn.xlf"
$ws.Range("E5").Value = $zhDatetime
$ws.Range("H5").Value = $epoch
$ws.Range("I5").Value = "Include"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdBase + $file1 + ".md", "", "", $file1 + ".md")
$ws.Hyperlinks.Add($ws.Range("B2"), $mdBase + $file1 + ".md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), $xlfZhBase + $file1 + "." + $hash1 + ".zh-cn.xlf", "", "", $file1 + "." + $hash1 + ".zh-cn.xlf")

$ws.Hyperlinks.Add($ws.Range("A3"), $mdBase + $file2 + ".md", "", "", $file2 + ".md")
$ws.Hyperlinks.Add($ws.Range("B3"), $mdBase + $file2 + ".md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), $xlfZhBase + $file2 + "." + $hash2 + ".zh-cn.xlf", "", "", $file2 + "." + $hash2 + ".zh-cn.xlf")

$ws.Hyperlinks.Add($ws.Range("A4"), $mdBase + $file3 + ".md", "", "", $file3 + ".md")
$ws.Hyperlinks.Add($ws.Range("B4"), $mdBase + $file3 + ".md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D4"), $xlfZhBase + $file3 + "." + $hash3 + ".zh-cn.xlf", "", "", $file3 + "." + $hash3 + ".zh-cn.xlf")

$ws.Hyperlinks.Add($ws.Range("A5"), $mdBase + $file4 + ".md", "", "", $file4 + ".md")
$ws.Hyperlinks.Add($ws.Range("B5"), $mdBase + $file4 + ".md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D5"), $xlfZhBase + $file4 + "." + $hash4 + ".zh-cn.xlf", "", "", $file4 + "." + $hash4 + ".zh-cn.xlf")

# =======================================================================
# Sheet "de-de"
# =======================================================================
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = $file1 + ".md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = $file1 + "." + $hash1 + ".de-de.xlf"
$ws.Range("E2").Value = $deDatetime
$ws.Range("H2").Value = $epoch
$ws.Range("I2").Value = "Include"

$ws.Range("A3").Value = $file2 + ".md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = $file2 + "." + $hash2 + ".de-de.xlf"
$ws.Range("E3").Value = $deDatetime
$ws.Range("H3").Value = $epoch
$ws.Range("I3").Value = "Include"

$ws.Range("A4").Value = $file3 + ".md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = $file3 + "." + $hash3 + ".de-de.xlf"
$ws.Range("E4").Value = $deDatetime
$ws.Range("H4").Value = $epoch
$ws.Range("I4").Value = "Include"

$ws.Range("A5").Value = $file4 + ".md"
$ws.Range("B5").Value = ".md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = $file4 + "." + $hash4 + ".de-de.xlf"
$ws.Range("E5").Value = $deDatetime
$ws.Range("H5").Value = $epoch
$ws.Range("I5").Value = "Include"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdBase + $file1 + ".md", "", "", $file1 + ".md")
$ws.Hyperlinks.Add($ws.Range("B2"), $mdBase + $file1 + ".md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), $xlfDeBase + $file1 + "." + $hash1 + ".de-de.xlf", "", "", $file1 + "." + $hash1 + ".de-de.xlf")

$ws.Hyperlinks.Add($ws.Range("A3"), $mdBase + $file2 + ".md", "", "", $file2 + ".md")
$ws.Hyperlinks.Add($ws.Range("B3"), $mdBase + $file2 + ".md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), $xlfDeBase + $file2 + "." + $hash2 + ".de-de.xlf", "", "", $file2 + "." + $hash2 + ".de-de.xlf")

$ws.Hyperlinks.Add($ws.Range("A4"), $mdBase + $file3 + ".md", "", "", $file3 + ".md")
$ws.Hyperlinks.Add($ws.Range("B4"), $mdBase + $file3 + ".md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D4"), $xlfDeBase + $file3 + "." + $hash3 + ".de-de.xlf", "", "", $file3 + "." + $hash3 + ".de-de.xlf")

$ws.Hyperlinks.Add($ws.Range("A5"), $mdBase + $file4 + ".md", "", "", $file4 + ".md")
$ws.Hyperlinks.Add($ws.Range("B5"), $mdBase + $file4 + ".md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D5"), $xlfDeBase + $file4 + "." + $hash4 + ".de-de.xlf", "", "", $file4 + "." + $hash4 + ".de-de.xlf")
